$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "mise en service d'une gateway neuve"
$ws.Range("D4").Value = "eviter reset intempestif"

$ws.Columns.Item(1).AutoFit() | Out-Null

$ws.Range("D4").Select() | Out-Null
